# Enhanced visualizations and improved model training
# Drop the four "raw DFT" columns (Total_Energy_Hartree, Solvation_Energy_eV,
# Surface_Area_A2, Molecular_Volume_A3) and shift the two absorption columns
# (Max_Absorption_nm, Max_f_osc) left into their place, then append a new
# derived column: Max_Excitation_eV.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns E:H entirely -- this shifts the old I (Max_Absorption_nm)
# and J (Max_f_osc) columns left to E and F, matching the diff.
$ws.Range("E1:H23").Delete()

# New third column with the excitation energy (eV), derived from
# Max_Absorption_nm via E = 1240 / lambda(nm), rounded to 2 dp.
# Match the bold/centered header formatting used by the other header cells
# (copy format from the neighboring F1 header, then set the new text).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Max_Excitation_eV"

$excitationEv = @(
    2.24,
    2.4,
    2.32,
    2.28,
    2.24,
    2.48,
    2.31,
    2.15,
    1.88,
    3.11,
    3.08,
    3.06,
    2.65,
    2.43,
    3.64,
    4.06,
    3.71,
    2.31,
    2.35,
    2.58,
    2.21,
    2.83
)

for ($i = 0; $i -lt $excitationEv.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $excitationEv[$i]
}
